# limpando o código e add funcoes
# Adds rows 9..25 to Sheet1, repeating the "Ar condicionado" product data
# that already exists in rows 3/5/6/7/8, and extends the used dimension
# accordingly (Excel recalculates the <dimension> automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$partNumber = "91697550"
$title = "Ar condicionado Split 24000 BTUs Quente e Frio 220V Series A1 TCL"

# Price text to use for each new row, keyed by row number (9..25)
$prices = @{
    9  = ",90"
    10 = ",90"
    11 = ",90"
    12 = ",90"
    13 = ",90"
    14 = ",90"
    15 = ",90"
    16 = ",90"
    17 = ",90"
    18 = ",90"
    19 = ",90"
    20 = ",90"
    21 = "3,949,90"
    22 = "3,949,90"
    23 = "3,949,90"
    24 = "3.949.90"
    25 = "3,949,90"
}

for ($row = 9; $row -le 25; $row++) {
    $ws.Cells.Item($row, 1).Value = "'" + $partNumber
    $ws.Cells.Item($row, 2).Value = $title
    $ws.Cells.Item($row, 3).Value = "'" + $prices[$row]
}
